$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: introduce "Initial Stocks Value" / "Remaining Stocks Value"
# in place of the old "Current Stocks" / "After Stocks", keeping "Total Sales" in C1.
$ws.Range("A1").Value = "Initial Stocks Value"
$ws.Range("B1").Value = "Remaining Stocks Value"
$ws.Range("C1").Value = "Total Sales"

# Row 4 headers: reorder columns so "Collectibles (Not Paid)" / "Due Collectibles (Deadline)"
# come first, followed by Store Name, Store Address, Order Received (Date), Sales Agent Name,
# Area Covered, Type of Outlet, Payment Method, and finally Discounts given per store.
$ws.Range("A4").Value = "Collectibles (Not Paid)"
$ws.Range("B4").Value = "Due Collectibles (Deadline)"
$ws.Range("C4").Value = "Store Name"
$ws.Range("D4").Value = "Store Address "
$ws.Range("E4").Value = "Order Received (Date)"
$ws.Range("F4").Value = "Sales Agent Name"
$ws.Range("G4").Value = "Area Covered"
$ws.Range("H4").Value = "Type of Outlet"
$ws.Range("I4").Value = "Payment Method"
$ws.Range("J4").Value = "Discounts given per store"

# Update the active selection to B1 (previously D6).
$ws.Range("B1").Select() | Out-Null
